$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.842.33"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.298.53"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'299.52"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'97.16"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "'35.94"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'6.75"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "2.656.03"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "2.295.66"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.777"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "42.827.17"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "'6.06"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").Value = "'67.91"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'242.62"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'25.12"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").Value = "'165.93"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'9.01"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'32.69"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'4.78"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("D36").Value = "'17.11"
$ws.Range("E36").Value = "  -4.36%  "
$ws.Range("D37").Value = "'2.38"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "'0.0686"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "2.007.40"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "'10.10"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  -5.20%  "
$ws.Range("D47").Value = "'17.18"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "2.525.28"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'52.99"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").Value = "'2.79"
$ws.Range("E51").Value = "  -7.20%  "
